$wb = $excel.ActiveWorkbook

# --- Sheet "results" ---
$wsResults = $wb.Worksheets.Item("results")
$wsResults.Range("B2").Value = 41.48399999999999
$wsResults.Range("C2").Value = 61.694
$wsResults.Range("D2").Value = 61.694
$wsResults.Range("E2").Value = 61.694
$wsResults.Range("F2").Value = 61.694
$wsResults.Range("G2").Value = 0.3897406228907532
$wsResults.Range("B3").Value = 32.656
$wsResults.Range("C3").Value = 47.382
$wsResults.Range("D3").Value = 47.382
$wsResults.Range("E3").Value = 47.382
$wsResults.Range("F3").Value = 47.382
$wsResults.Range("G3").Value = 0.3607545320921117
$wsResults.Range("B4").Value = 33.48399999999999
$wsResults.Range("C4").Value = 41.382
$wsResults.Range("D4").Value = 41.382
$wsResults.Range("E4").Value = 41.382
$wsResults.Range("F4").Value = 41.382
$wsResults.Range("G4").Value = 0.1886990801576875
$wsResults.Range("B5").Value = 50.038
$wsResults.Range("C5").Value = 48.20999999999999
$wsResults.Range("D5").Value = 48.20999999999999
$wsResults.Range("E5").Value = 48.20999999999999
$wsResults.Range("F5").Value = 48.20999999999999
$wsResults.Range("G5").Value = 0.02922578840081569
$wsResults.Range("B6").Value = 41.694
$wsResults.Range("C6").Value = 45.452
$wsResults.Range("D6").Value = 45.452
$wsResults.Range("E6").Value = 45.452
$wsResults.Range("F6").Value = 45.452
$wsResults.Range("G6").Value = 0.07210629826833603

# --- Sheet "stats" ---
$wsStats = $wb.Worksheets.Item("stats")
$wsStats.Range("C2").Value = 103
$wsStats.Range("D2").Value = 0.008189005544409156
$wsStats.Range("E2").Value = 0.105924578034319
$wsStats.Range("F2").Value = 103
$wsStats.Range("G2").Value = 0.005882885656319559
$wsStats.Range("H2").Value = 0.06969061761628836
$wsStats.Range("I2").Value = 0.005414729355834424
$wsStats.Range("J2").Value = 0.01658011670224369
$wsStats.Range("K2").Value = 0.002096734941005707
$wsStats.Range("C3").Value = 841
$wsStats.Range("D3").Value = 0.001355051877908409
$wsStats.Range("E3").Value = 0.4458491399418563
$wsStats.Range("F3").Value = 841
$wsStats.Range("G3").Value = 0.03444239869713783
$wsStats.Range("H3").Value = 0.2877437327988446
$wsStats.Range("I3").Value = 0.0161881810054183
$wsStats.Range("J3").Value = 0.0571952520404011
$wsStats.Range("K3").Value = 0.01244004664476961
$wsStats.Range("C4").Value = 103
$wsStats.Range("D4").Value = 0.005695371772162616
$wsStats.Range("E4").Value = 0.1002145250095055
$wsStats.Range("F4").Value = 103
$wsStats.Range("G4").Value = 0.01339302223641425
$wsStats.Range("H4").Value = 0.06159640965051949
$wsStats.Range("I4").Value = 0.001973273698240519
$wsStats.Range("J4").Value = 0.01416924560908228
$wsStats.Range("K4").Value = 0.002171794883906841
$wsStats.Range("C5").Value = 841
$wsStats.Range("D5").Value = 0.001980427536182106
$wsStats.Range("E5").Value = 0.4499561260454357
$wsStats.Range("F5").Value = 841
$wsStats.Range("G5").Value = 0.03585337020922452
$wsStats.Range("H5").Value = 0.2926456088898703
$wsStats.Range("I5").Value = 0.01467407483141869
$wsStats.Range("J5").Value = 0.05429538327734917
$wsStats.Range("K5").Value = 0.01307916862424463
$wsStats.Range("E6").Value = 1.062308883643709
$wsStats.Range("C7").Value = 59
$wsStats.Range("D7").Value = 0.003549611894413829
$wsStats.Range("E7").Value = 0.0703013630118221
$wsStats.Range("F7").Value = 59
$wsStats.Range("G7").Value = 0.003958272864110768
$wsStats.Range("H7").Value = 0.05023069749586284
$wsStats.Range("I7").Value = 0.001357558649033308
$wsStats.Range("J7").Value = 0.009682819480076432
$wsStats.Range("K7").Value = 0.001297194859944284
$wsStats.Range("C8").Value = 684
$wsStats.Range("D8").Value = 0.001199833932332695
$wsStats.Range("E8").Value = 0.3780388489831239
$wsStats.Range("F8").Value = 684
$wsStats.Range("G8").Value = 0.02844088862184435
$wsStats.Range("H8").Value = 0.245931152603589
$wsStats.Range("I8").Value = 0.01004169916268438
$wsStats.Range("J8").Value = 0.05189107090700418
$wsStats.Range("K8").Value = 0.01022336073219776
$wsStats.Range("C9").Value = 59
$wsStats.Range("D9").Value = 0.00385952671058476
$wsStats.Range("E9").Value = 0.05054009100422263
$wsStats.Range("F9").Value = 59
$wsStats.Range("G9").Value = 0.003202077932655811
$wsStats.Range("H9").Value = 0.03270551166497171
$wsStats.Range("I9").Value = 0.001332056242972612
$wsStats.Range("J9").Value = 0.008711365284398198
$wsStats.Range("K9").Value = 0.00121563533321023
$wsStats.Range("C10").Value = 684
$wsStats.Range("D10").Value = 0.001812252099625766
$wsStats.Range("E10").Value = 0.4078752159839496
$wsStats.Range("F10").Value = 684
$wsStats.Range("G10").Value = 0.03220666048582643
$wsStats.Range("H10").Value = 0.2666687830351293
$wsStats.Range("I10").Value = 0.0129236081847921
$wsStats.Range("J10").Value = 0.04882711649406701
$wsStats.Range("K10").Value = 0.01174747699405998
$wsStats.Range("E11").Value = 0.8592982107074931
$wsStats.Range("C12").Value = 85
$wsStats.Range("D12").Value = 0.002475501387380064
$wsStats.Range("E12").Value = 0.04337085306178778
$wsStats.Range("F12").Value = 85
$wsStats.Range("G12").Value = 0.002991055836901069
$wsStats.Range("H12").Value = 0.02836394682526588
$wsStats.Range("I12").Value = 0.0009896616684272885
$wsStats.Range("J12").Value = 0.006671224138699472
$wsStats.Range("K12").Value = 0.001111330115236342
$wsStats.Range("C13").Value = 617
$wsStats.Range("D13").Value = 0.0008665162604302168
$wsStats.Range("E13").Value = 0.2870514150708914
$wsStats.Range("F13").Value = 617
$wsStats.Range("G13").Value = 0.02220991323702037
$wsStats.Range("H13").Value = 0.1892062119441107
$wsStats.Range("I13").Value = 0.007630242151208222
$wsStats.Range("J13").Value = 0.03455374017357826
$wsStats.Range("K13").Value = 0.008201737073250115
$wsStats.Range("C14").Value = 85
$wsStats.Range("D14").Value = 0.002971776877529919
$wsStats.Range("E14").Value = 0.04787431401200593
$wsStats.Range("F14").Value = 85
$wsStats.Range("G14").Value = 0.003168480703607202
$wsStats.Range("H14").Value = 0.03155661758501083
$wsStats.Range("I14").Value = 0.001137028448283672
$wsStats.Range("J14").Value = 0.007419725530780852
$wsStats.Range("K14").Value = 0.001158107188530266
$wsStats.Range("C15").Value = 617
$wsStats.Range("D15").Value = 0.001197399804368615
$wsStats.Range("E15").Value = 0.2643401949899271
$wsStats.Range("F15").Value = 617
$wsStats.Range("G15").Value = 0.02032285393215716
$wsStats.Range("H15").Value = 0.1737123600905761
$wsStats.Range("I15").Value = 0.00831927452236414
$wsStats.Range("J15").Value = 0.03230136178899556
$wsStats.Range("K15").Value = 0.007420989568345249
$wsStats.Range("E16").Value = 1.176353588118218
$wsStats.Range("C17").Value = 97
$wsStats.Range("D17").Value = 0.002722627250477672
$wsStats.Range("E17").Value = 0.04655965196434408
$wsStats.Range("F17").Value = 97
$wsStats.Range("G17").Value = 0.00305652036331594
$wsStats.Range("H17").Value = 0.03000462322961539
$wsStats.Range("I17").Value = 0.001027111429721117
$wsStats.Range("J17").Value = 0.007854246068745852
$wsStats.Range("K17").Value = 0.001157607650384307
$wsStats.Range("C18").Value = 691
$wsStats.Range("D18").Value = 0.0008000775706022978
$wsStats.Range("E18").Value = 0.2823822669452056
$wsStats.Range("F18").Value = 691
$wsStats.Range("G18").Value = 0.02191348955966532
$wsStats.Range("H18").Value = 0.1842833496630192
$wsStats.Range("I18").Value = 0.008032356854528189
$wsStats.Range("J18").Value = 0.03575757960788906
$wsStats.Range("K18").Value = 0.00800105242524296
$wsStats.Range("C19").Value = 97
$wsStats.Range("D19").Value = 0.003338298061862588
$wsStats.Range("E19").Value = 0.06098816101439297
$wsStats.Range("F19").Value = 97
$wsStats.Range("G19").Value = 0.003850886365398765
$wsStats.Range("H19").Value = 0.03933807380963117
$wsStats.Range("I19").Value = 0.001321790041401982
$wsStats.Range("J19").Value = 0.01039966521784663
$wsStats.Range("K19").Value = 0.00158635200932622
$wsStats.Range("C20").Value = 691
$wsStats.Range("D20").Value = 0.001357300905510783
$wsStats.Range("E20").Value = 0.2997041610069573
$wsStats.Range("F20").Value = 691
$wsStats.Range("G20").Value = 0.02359679853543639
$wsStats.Range("H20").Value = 0.1951985992491245
$wsStats.Range("I20").Value = 0.009640221367590129
$wsStats.Range("J20").Value = 0.03715040034148842
$wsStats.Range("K20").Value = 0.008432312519289553
$wsStats.Range("E21").Value = 1.37943260197062
$wsStats.Range("C22").Value = 80
$wsStats.Range("D22").Value = 0.002464421442709863
$wsStats.Range("E22").Value = 0.03949186997488141
$wsStats.Range("F22").Value = 80
$wsStats.Range("G22").Value = 0.002676427946425974
$wsStats.Range("H22").Value = 0.0256324663059786
$wsStats.Range("I22").Value = 0.0009342546109110117
$wsStats.Range("J22").Value = 0.006440620636567473
$wsStats.Range("K22").Value = 0.0009781229309737682
$wsStats.Range("C23").Value = 645
$wsStats.Range("D23").Value = 0.0008245084900408983
$wsStats.Range("E23").Value = 0.2754489569924772
$wsStats.Range("F23").Value = 645
$wsStats.Range("G23").Value = 0.02135876473039389
$wsStats.Range("H23").Value = 0.1799972007283941
$wsStats.Range("I23").Value = 0.008024295908398926
$wsStats.Range("J23").Value = 0.0347382293548435
$wsStats.Range("K23").Value = 0.007764300680719316
$wsStats.Range("C24").Value = 80
$wsStats.Range("D24").Value = 0.002544151502661407
$wsStats.Range("E24").Value = 0.03929092094767839
$wsStats.Range("F24").Value = 80
$wsStats.Range("G24").Value = 0.002686112653464079
$wsStats.Range("H24").Value = 0.02556078461930156
$wsStats.Range("I24").Value = 0.0009594347793608904
$wsStats.Range("J24").Value = 0.006331092794425786
$wsStats.Range("K24").Value = 0.000925549422390759
$wsStats.Range("C25").Value = 645
$wsStats.Range("D25").Value = 0.001189700677059591
$wsStats.Range("E25").Value = 0.2784969829954207
$wsStats.Range("F25").Value = 645
$wsStats.Range("G25").Value = 0.02145504322834313
$wsStats.Range("H25").Value = 0.1806719757150859
$wsStats.Range("I25").Value = 0.009620036697015166
$wsStats.Range("J25").Value = 0.03499171510338783
$wsStats.Range("K25").Value = 0.007946688449010253
$wsStats.Range("E26").Value = 1.507189766038209
$wsStats.Range("C27").Value = 84.8
$wsStats.Range("D27").Value = 0.003880233503878117
$wsStats.Range("E27").Value = 0.06112966320943088
$wsStats.Range("F27").Value = 84.8
$wsStats.Range("G27").Value = 0.003713032533414662
$wsStats.Range("H27").Value = 0.04078447029460221
$wsStats.Range("I27").Value = 0.00194466314278543
$wsStats.Range("J27").Value = 0.009445805405266582
$wsStats.Range("K27").Value = 0.001328198099508882
$wsStats.Range("C28").Value = 695.6
$wsStats.Range("D28").Value = 0.001009197626262903
$wsStats.Range("E28").Value = 0.3337541255867109
$wsStats.Range("F28").Value = 695.6
$wsStats.Range("G28").Value = 0.02567309096921235
$wsStats.Range("H28").Value = 0.2174323295475915
$wsStats.Range("I28").Value = 0.009983355016447603
$wsStats.Range("J28").Value = 0.04282717441674322
$wsStats.Range("K28").Value = 0.009326099511235952
$wsStats.Range("C29").Value = 84.8
$wsStats.Range("D29").Value = 0.003681824984960258
$wsStats.Range("E29").Value = 0.05978160239756107
$wsStats.Range("F29").Value = 84.8
$wsStats.Range("G29").Value = 0.005260115978308022
$wsStats.Range("H29").Value = 0.03815147946588695
$wsStats.Range("I29").Value = 0.001344716642051935
$wsStats.Range("J29").Value = 0.00940621888730675
$wsStats.Range("K29").Value = 0.001411487767472863
$wsStats.Range("C30").Value = 695.6
$wsStats.Range("D30").Value = 0.001507416204549372
$wsStats.Range("E30").Value = 0.3400745362043381
$wsStats.Range("F30").Value = 695.6
$wsStats.Range("G30").Value = 0.02668694527819753
$wsStats.Range("H30").Value = 0.2217794653959572
$wsStats.Range("I30").Value = 0.01103544312063605
$wsStats.Range("J30").Value = 0.0415131954010576
$wsStats.Range("K30").Value = 0.009725327230989933
$wsStats.Range("E31").Value = 1.19691661009565

